# Auto-generated edit script applying the Marilith_Profits.xlsx diff.
# Updates/adds/removes specific H/I/J/K/L/M/N numeric cells across the 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3542.077
$ws.Range("J40").Value = 3587.3333
$ws.Range("L40").Value = 3587.3333
$ws.Range("N40").Value = -3937.3333
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H100").Value = 19999.5
$ws.Range("J100").Value = 30000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -31082
$ws.Range("H112").Value = 2138.889
$ws.Range("J112").Value = 2321.4285
$ws.Range("L112").Value = 6964.2855
$ws.Range("N112").Value = -9180.2855
$ws.Range("H113").Value = 4115.8887
$ws.Range("J113").Value = 5499.3335
$ws.Range("L113").Value = 5499.3335
$ws.Range("N113").Value = -12007.3335
$ws.Range("H125").Value = 11333
$ws.Range("I125").Value = 10000
$ws.Range("K125").Value = 90000
$ws.Range("M125").Value = -87540
$ws.Range("H131").Value = 2999.8333
$ws.Range("I131").Value = 2559.8
$ws.Range("J131").Value = 5200
$ws.Range("K131").Value = 7679.400000000001
$ws.Range("L131").Value = 15600
$ws.Range("M131").Value = -2639.400000000001
$ws.Range("N131").Value = -25680
$ws.Range("H137").Value = 2165.1333
$ws.Range("I137").Value = 1723.1666
$ws.Range("K137").Value = 5169.4998
$ws.Range("M137").Value = -2619.4998
$ws.Range("H138").Value = 3200.0667
$ws.Range("J138").Value = 3461.7693
$ws.Range("L138").Value = 10385.3079
$ws.Range("N138").Value = -20665.3079

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8502.654
$ws.Range("I32").Value = 8502.654
$ws.Range("K32").Value = 8502.654
$ws.Range("M32").Value = -8215.654
$ws.Range("H45").Value = 2922.625
$ws.Range("I45").Value = 2731
$ws.Range("K45").Value = 2731
$ws.Range("M45").Value = -2354
$ws.Range("H74").Value = 1920
$ws.Range("I74").Value = 1700.125
$ws.Range("J74").Value = 2799.5
$ws.Range("K74").Value = 1700.125
$ws.Range("L74").Value = 2799.5
$ws.Range("M74").Value = -826.125
$ws.Range("N74").Value = -4547.5
$ws.Range("H77").Value = 1920
$ws.Range("I77").Value = 1700.125
$ws.Range("J77").Value = 2799.5
$ws.Range("K77").Value = 8500.625
$ws.Range("L77").Value = 13997.5
$ws.Range("M77").Value = -4132.625
$ws.Range("N77").Value = -22733.5
$ws.Range("H97").Value = 1224.2858
$ws.Range("J97").Value = 2945
$ws.Range("L97").Value = 2945
$ws.Range("N97").Value = -3937

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2010.95
$ws.Range("I86").Value = 2233.7856
$ws.Range("J86").Value = 1491
$ws.Range("K86").Value = 2233.7856
$ws.Range("L86").Value = 1491
$ws.Range("M86").Value = -1110.7856
$ws.Range("N86").Value = -3737
$ws.Range("H89").Value = 2010.95
$ws.Range("I89").Value = 2233.7856
$ws.Range("J89").Value = 1491
$ws.Range("K89").Value = 11168.928
$ws.Range("L89").Value = 7455
$ws.Range("M89").Value = -5552.928
$ws.Range("N89").Value = -18687

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4782.25
$ws.Range("J31").Value = 5921.6665
$ws.Range("L31").Value = 5921.6665
$ws.Range("N31").Value = -6511.6665
$ws.Range("H34").Value = 4782.25
$ws.Range("J34").Value = 5921.6665
$ws.Range("L34").Value = 5921.6665
$ws.Range("N34").Value = -6325.6665
$ws.Range("H99").Value = 2650
$ws.Range("I99").Value = 2250
$ws.Range("J99").Value = 2783.3333
$ws.Range("K99").Value = 2250
$ws.Range("L99").Value = 2783.3333
$ws.Range("M99").Value = -752
$ws.Range("N99").Value = -5779.3333
$ws.Range("H126").Value = 2650
$ws.Range("I126").Value = 2250
$ws.Range("J126").Value = 2783.3333
$ws.Range("K126").Value = 6750
$ws.Range("L126").Value = 8349.999899999999
$ws.Range("M126").Value = -4280
$ws.Range("N126").Value = -13289.9999
$ws.Range("H132").Value = 2881.5
$ws.Range("I132").Value = 2881.5
$ws.Range("K132").Value = 8644.5
$ws.Range("M132").Value = -6114.5
$ws.Range("H134").Value = 2026.579
$ws.Range("I134").Value = 2176.5
$ws.Range("J134").Value = 1227
$ws.Range("K134").Value = 6529.5
$ws.Range("L134").Value = 3681
$ws.Range("M134").Value = -3994.5
$ws.Range("N134").Value = -8751

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 997.5
$ws.Range("I120").Value = 997.5
$ws.Range("K120").Value = 2992.5
$ws.Range("M120").Value = 1845.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -712
$ws.Range("N19").ClearContents()
$ws.Range("H122").Value = 41588.332
$ws.Range("I122").Value = 40572.734
$ws.Range("K122").Value = 121718.202
$ws.Range("M122").Value = -119268.202
$ws.Range("H126").Value = 1400
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 4200
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -9140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1549.1666
$ws.Range("I7").Value = 1858
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 1858
$ws.Range("L7").Value = 5
$ws.Range("M7").Value = -1746
$ws.Range("N7").Value = -229
$ws.Range("H16").Value = 2309.75
$ws.Range("I16").Value = 2635.889
$ws.Range("K16").Value = 2635.889
$ws.Range("M16").Value = -2465.889
$ws.Range("H22").Value = 1212.95
$ws.Range("I22").Value = 1078.7273
$ws.Range("J22").Value = 1377
$ws.Range("K22").Value = 1078.7273
$ws.Range("L22").Value = 1377
$ws.Range("M22").Value = -783.7273
$ws.Range("N22").Value = -1967
$ws.Range("H27").Value = 1212.95
$ws.Range("I27").Value = 1078.7273
$ws.Range("J27").Value = 1377
$ws.Range("K27").Value = 1078.7273
$ws.Range("L27").Value = 1377
$ws.Range("M27").Value = -971.7273
$ws.Range("N27").Value = -1591
$ws.Range("H40").Value = 4400.7334
$ws.Range("I40").Value = 3001
$ws.Range("K40").Value = 3001
$ws.Range("M40").Value = -2865
$ws.Range("H46").Value = 2149.125
$ws.Range("I46").Value = 1097.5
$ws.Range("J46").Value = 2499.6667
$ws.Range("K46").Value = 1097.5
$ws.Range("L46").Value = 2499.6667
$ws.Range("M46").Value = -909.5
$ws.Range("N46").Value = -2875.6667
$ws.Range("H55").Value = 745.8570999999999
$ws.Range("I55").Value = 722
$ws.Range("K55").Value = 722
$ws.Range("M55").Value = -549
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 1549.1666
$ws.Range("I126").Value = 1858
$ws.Range("J126").Value = 5
$ws.Range("K126").Value = 5574
$ws.Range("L126").Value = 15
$ws.Range("M126").Value = -3104
$ws.Range("N126").Value = -4955
$ws.Range("H136").Value = 3517.25
$ws.Range("I136").Value = 3517.25
$ws.Range("K136").Value = 10551.75
$ws.Range("M136").Value = -8001.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 9300
$ws.Range("J43").Value = 9300
$ws.Range("L43").Value = 9300
$ws.Range("N43").Value = -9598
$ws.Range("H126").Value = 2574.25
$ws.Range("I126").Value = 2599
$ws.Range("K126").Value = 7797
$ws.Range("M126").Value = -5327
$ws.Range("H132").Value = 500
$ws.Range("I132").Value = 500
$ws.Range("K132").Value = 1500
$ws.Range("M132").Value = 1030
$ws.Range("H136").Value = 3945.5
$ws.Range("I136").Value = 3945.5
$ws.Range("K136").Value = 11836.5
$ws.Range("M136").Value = -9286.5

Write-Output "Applied all Marilith_Profits.xlsx cell updates"